$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# Order chosen so the rebuilt shared-string table comes out in the same
# order as the target workbook (Comedia, El show de Truman, Thriller,
# Padre no hay más que uno, padre_uno url, show_truman url).
$ws.Range("B2").Value = "Comedia"
$ws.Range("A1").Value = "El show de Truman"
$ws.Range("B1").Value = "Thriller"
$ws.Range("A2").Value = "Padre no hay más que uno"
$ws.Range("C2").Value = "https://github.com/israel-android/excel_pruebas/blob/main/padre_uno.jpg"
$ws.Range("C1").Value = "https://github.com/israel-android/excel_pruebas/blob/main/show_truman.jpg"

# --- Hyperlinks (also applies the built-in Hyperlink cell style) -----
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/israel-android/excel_pruebas/blob/main/padre_uno.jpg")
$ws.Hyperlinks.Add($ws.Range("C1"), "https://github.com/israel-android/excel_pruebas/blob/main/show_truman.jpg")

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.16666666666667
$ws.Columns.Item(2).ColumnWidth = 29.33333333333333
$ws.Columns.Item(3).ColumnWidth = 72.16666666666667

# --- Selection ---------------------------------------------------------
$ws.Range("C6").Select() | Out-Null
